# Error Calculations and Plots
# This "missing_data" worksheet simulates randomly-removed observations.
# Two rows (RM 232 and SC 92) are removed from the dataset, and the set of
# cells considered "missing" in column A (header label) is updated for the
# remaining rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RM 232" row (row 26). Everything below shifts up by one.
$ws.Rows.Item(26).Delete()

# After the first deletion, the former "SC 92" row is now row 27. Remove it
# too; everything below shifts up by one more.
$ws.Rows.Item(27).Delete()

# Update which cells are treated as "missing" (blank) in column B for the
# remaining rows.
$ws.Range("B26").Value = ""          # SC 5  -> now missing
$ws.Range("B27").Value = -20.4       # SC 101 -> now present
$ws.Range("B29").Value = ""          # SC 119 -> now missing
